# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header style used by the other columns (e.g. G1's bold/bordered style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's formatting (style) onto H1, then set text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Values for H2:H28, in row order (row 2 is the first data row).
$saveValues = @(0,0,1,1,1,0,0,0,0,1,0,0,1,0,1,0,0,0,0,0,1,0,0,0,0,1,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
